$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D7").Value = 0
$ws.Range("J7").Value = 612
$ws.Range("B8").Value = 372.2
$ws.Range("D8").Value = 0
$ws.Range("F8").Value = 372.2
$ws.Range("J8").Value = 372.2
$ws.Range("K8").Value = 1
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 372.2
$ws.Range("P8").Value = 0
$ws.Range("B9").Value = 372.2
$ws.Range("D9").Value = 0
$ws.Range("F9").Value = 372.2
$ws.Range("J9").Value = 372.2
$ws.Range("K9").Value = 1
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 372.2
$ws.Range("P9").Value = 0
$ws.Range("B10").Value = 244
$ws.Range("C10").Value = 1
$ws.Range("E10").Value = 0
$ws.Range("J10").Value = 244
$ws.Range("K10").Value = 1
$ws.Range("M10").Value = 0
$ws.Range("B11").Value = 243.5
$ws.Range("C11").Value = 1
$ws.Range("E11").Value = 0
$ws.Range("J11").Value = 243.5
$ws.Range("K11").Value = 1
$ws.Range("M11").Value = 0
$ws.Range("B12").Value = 149.8
$ws.Range("C12").Value = 1
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 149.8
$ws.Range("G12").Value = 1
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 149.8
$ws.Range("K12").Value = 1
$ws.Range("M12").Value = 0
$ws.Range("B13").Value = 144.4
$ws.Range("C13").Value = 1
$ws.Range("E13").Value = 0
$ws.Range("J13").Value = 144.4
$ws.Range("K13").Value = 1
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 144.4
$ws.Range("O13").Value = 1
$ws.Range("Q13").Value = 0
$ws.Range("B14").Value = 82
$ws.Range("C14").Value = 1
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 82
$ws.Range("G14").Value = 1
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 82
$ws.Range("K14").Value = 1
$ws.Range("M14").Value = 0
$ws.Range("AC14").Value = 0
$ws.Range("B15").Value = 80
$ws.Range("C15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("J15").Value = 80
$ws.Range("K15").Value = 1
$ws.Range("M15").Value = 0
$ws.Range("AC15").Value = 0
$ws.Range("B16").Value = 47.9
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 47.9
$ws.Range("G16").Value = 1
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 47.9
$ws.Range("K16").Value = 1
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 47.9
$ws.Range("O16").Value = 1
$ws.Range("Q16").Value = 0
$ws.Range("B17").Value = 47.5
$ws.Range("C17").Value = 1
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 47.5
$ws.Range("G17").Value = 1
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 47.5
$ws.Range("K17").Value = 1
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 47.5
$ws.Range("O17").Value = 1
$ws.Range("Q17").Value = 0
$ws.Range("U25").Value = 3
$ws.Range("Y25").Value = 1
$ws.Range("U32").Value = 11
$ws.Range("B33").Value = 664.7415060000021
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").Value = 10
$ws.Range("B34").Value = 0
$ws.Range("C34").Value = 0
$ws.Range("E34").Value = 11
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").Value = 10
$ws.Range("R35").Value = 311.9558659999997
$ws.Range("S35").Value = 1
$ws.Range("U35").Value = 0
$ws.Range("D36").Value = 19
$ws.Range("F36").Value = 110.96
$ws.Range("B37").Value = 115.5
$ws.Range("F37").Value = 301.8888139999999
$ws.Range("I37").Value = 2
$ws.Range("J37").Value = 103.1
$ws.Range("N37").Value = 249.7482900305193
$ws.Range("U37").Value = 1
$ws.Range("B41").Value = 0
$ws.Range("C41").Value = 0
$ws.Range("E41").Value = 3
$ws.Range("F41").Value = 54.18
$ws.Range("J41").Value = 54.18
$ws.Range("N41").Value = 54.17999999999995
$ws.Range("V41").Value = 270.4792560000037
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("N46").Value = 0
$ws.Range("O46").Value = 0
$ws.Range("U47").Value = 6
$ws.Range("AC47").Value = 7
$ws.Range("Q48").Value = 4
$ws.Range("Y48").Value = 3
$ws.Range("AC48").Value = 4
$ws.Range("J49").Value = 218.6226880000033
$ws.Range("M50").Value = 3
$ws.Range("Q50").Value = 5
$ws.Range("R50").Value = 0
$ws.Range("S50").Value = 0
$ws.Range("U50").Value = 5
$ws.Range("U51").Value = 11
$ws.Range("Z51").Value = 30.94000000000102
$ws.Range("Q52").Value = 5
$ws.Range("R52").Value = 0
$ws.Range("S52").Value = 0
$ws.Range("U52").Value = 5
$ws.Range("Q54").Value = 4
$ws.Range("Y54").Value = 2
$ws.Range("M55").Value = 1
$ws.Range("Y55").Value = 1
$ws.Range("Q56").Value = 2
$ws.Range("Y56").Value = 2
$ws.Range("E58").Value = 5
$ws.Range("M58").Value = 4
$ws.Range("Q58").Value = 5
$ws.Range("Y58").Value = 6
$ws.Range("AC58").Value = 6
$ws.Range("B59").Value = 0
$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("E59").Value = 3
